$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115
$ws.Range("D115").Value = 45021
$ws.Range("L115").Value = "Especial"
$ws.Range("M115").Value = 16
$ws.Range("N115").Value = 230000
$ws.Range("O115").Value = 240000
$ws.Range("P115").Value = 235000
$ws.Range("R115").Value = "Región de O'Higgins"
$ws.Range("S115").Value = 522

# Row 116
$ws.Range("D116").Value = 45021
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 20
$ws.Range("N116").Value = 200000
$ws.Range("O116").Value = 210000
$ws.Range("P116").Value = 205000
$ws.Range("R116").Value = "Región de O'Higgins"
$ws.Range("S116").Value = 456

# Row 117
$ws.Range("D117").Value = 45021
$ws.Range("K117").Value = "Angeleno"
$ws.Range("L117").Value = "Segunda"
$ws.Range("M117").Value = 16
$ws.Range("N117").Value = 170000
$ws.Range("O117").Value = 180000
$ws.Range("P117").Value = 175000
$ws.Range("S117").Value = 389

# Row 118
$ws.Range("D118").Value = 45015
$ws.Range("K118").Value = "Angeleno"
$ws.Range("L118").Value = "Primera"
$ws.Range("M118").Value = 20
$ws.Range("N118").Value = 220000
$ws.Range("O118").Value = 230000
$ws.Range("P118").Value = 225000
$ws.Range("R118").Value = "Región Metropolitana"
$ws.Range("S118").Value = 500

# Row 119
$ws.Range("D119").Value = 45015
$ws.Range("L119").Value = "Segunda"
$ws.Range("M119").Value = 14
$ws.Range("N119").Value = 180000
$ws.Range("O119").Value = 190000
$ws.Range("P119").Value = 185000
$ws.Range("R119").Value = "Región Metropolitana"
$ws.Range("S119").Value = 411

# Row 120
$ws.Range("A120").Value = 2
$ws.Range("B120").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 44952
$ws.Range("E120").Value = 4
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100103
$ws.Range("H120").Value = "Frutos de hueso (carozo)"
$ws.Range("I120").Value = 100103002
$ws.Range("J120").Value = "Ciruela"
$ws.Range("K120").Value = "Black Amber"
$ws.Range("L120").Value = "Primera"
$ws.Range("M120").Value = 10
$ws.Range("N120").Value = 300000
$ws.Range("O120").Value = 310000
$ws.Range("P120").Value = 305000
$ws.Range("Q120").Value = "`$/bins (450 kilos)"
$ws.Range("R120").Value = "Región de O'Higgins"
$ws.Range("S120").Value = 678
$ws.Range("T120").Value = 450

# Row 121
$ws.Range("A121").Value = 2
$ws.Range("B121").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C121").Value = "Coquimbo"
$ws.Range("D121").Value = 44952
$ws.Range("D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E121").Value = 4
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100103
$ws.Range("H121").Value = "Frutos de hueso (carozo)"
$ws.Range("I121").Value = 100103002
$ws.Range("J121").Value = "Ciruela"
$ws.Range("K121").Value = "Black Amber"
$ws.Range("L121").Value = "Segunda"
$ws.Range("M121").Value = 10
$ws.Range("N121").Value = 230000
$ws.Range("O121").Value = 240000
$ws.Range("P121").Value = 235000
$ws.Range("Q121").Value = "`$/bins (450 kilos)"
$ws.Range("R121").Value = "Región de O'Higgins"
$ws.Range("S121").Value = 522
$ws.Range("T121").Value = 450

# Row 122
$ws.Range("A122").Value = 2
$ws.Range("B122").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 45007
$ws.Range("D122").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E122").Value = 4
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100103
$ws.Range("H122").Value = "Frutos de hueso (carozo)"
$ws.Range("I122").Value = 100103002
$ws.Range("J122").Value = "Ciruela"
$ws.Range("K122").Value = "Angeleno"
$ws.Range("L122").Value = "Primera"
$ws.Range("M122").Value = 18
$ws.Range("N122").Value = 200000
$ws.Range("O122").Value = 210000
$ws.Range("P122").Value = 205000
$ws.Range("Q122").Value = "`$/bins (450 kilos)"
$ws.Range("R122").Value = "Región de O'Higgins"
$ws.Range("S122").Value = 456
$ws.Range("T122").Value = 450

# Row 123
$ws.Range("A123").Value = 2
$ws.Range("B123").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C123").Value = "Coquimbo"
$ws.Range("D123").Value = 45007
$ws.Range("D123").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E123").Value = 4
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100103
$ws.Range("H123").Value = "Frutos de hueso (carozo)"
$ws.Range("I123").Value = 100103002
$ws.Range("J123").Value = "Ciruela"
$ws.Range("K123").Value = "Angeleno"
$ws.Range("L123").Value = "Segunda"
$ws.Range("M123").Value = 14
$ws.Range("N123").Value = 180000
$ws.Range("O123").Value = 190000
$ws.Range("P123").Value = 185000
$ws.Range("Q123").Value = "`$/bins (450 kilos)"
$ws.Range("R123").Value = "Región de O'Higgins"
$ws.Range("S123").Value = 411
$ws.Range("T123").Value = 450
